$wb = $excel.ActiveWorkbook

# --- Rename sheets ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws1.Name = "copy editor"
$ws2.Name = "Chris"

# --- Populate "Chris" sheet (formerly Sheet2) with review notes -------
$xlRight = -4152
$xlCenter = -4108

function Set-Row($r, $aText, $bNum, $cText) {
    $ws2.Range("A$r").HorizontalAlignment = $xlRight
    $ws2.Range("B$r").HorizontalAlignment = $xlCenter
    if ($aText -ne $null) {
        $ws2.Range("A$r").Value = $aText
    }
    if ($bNum -ne $null) {
        $ws2.Range("B$r").Value = $bNum
    }
    if ($cText -ne $null) {
        $ws2.Range("C$r").Value = $cText
    }
}

# Row 1 is the header: "Page" label lives in column B (centered), not C.
$ws2.Range("A1").HorizontalAlignment = $xlRight
$ws2.Range("B1").HorizontalAlignment = $xlCenter
$ws2.Range("B1").Value = "Page"

Set-Row 2 $null 14 "pre….. ? Risk …s not YLDs"
Set-Row 3 "x" 15 "diverse"
Set-Row 4 $null 16 "use of word novel"
Set-Row 5 $null 17 "citation for YLD paper"
Set-Row 6 $null 17 "many …. TV …. This"
Set-Row 7 $null 17 "… data are not for these years"
Set-Row 8 $null 17 "define pure number"
Set-Row 9 "x" 18 "figure"
Set-Row 10 $null 18 "why do we allow for 1990 data … 1980 but not for 2010"
Set-Row 11 "x" 19 "country level fixed effects"
Set-Row 12 "x" 20 "figure scale"
Set-Row 13 "x" 21 "wording strange ""might be"""
Set-Row 14 $null 21 "why give statistics of ACS and AMS?"
Set-Row 15 $null 24 "not correct, prosanto made this. DisMod1. this is not correct."
Set-Row 16 "x" 24 "Lotus ref"
Set-Row 17 "x" 25 "number of national burden of disease studies"

# --- Move active tab / selection to the "Chris" sheet -----------------
$ws2.Activate()
$ws2.Range("L19").Select()
